$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.298.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4726"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.57%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2873"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06478"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07785"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.864.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7173"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.125"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "280.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.286.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007465"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.107.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.97%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.240"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.251"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.971"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.878"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09635"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.77%  "

$ws.Range("E30").Value = "  -2.50%  "

$ws.Range("E31").Value = "  -1.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.207"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.113"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04780"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.116"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6834"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.11%  "

$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("E38").Value = "  -0.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.841"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.224"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.929"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4190"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9988"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8246"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.579"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.961"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05767"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "881.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.88%  "
